$wb = $excel.ActiveWorkbook

# --- Rename the second sheet tab ---
$wsInclude = $wb.Worksheets.Item("Include from unknown")
$wsInclude.Name = "Include #0"

# --- Update Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version 1.0.0 -> 1.0.1
$wsMeta.Range("B3").Value = "1.0.1"

# Contact value change
$wsMeta.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# Insert a new row after row 10 (Contact) for Jurisdiction
$wsMeta.Rows.Item(11).Insert()

# Copy style from the row above (Contact row) for the new Jurisdiction row cells
$wsMeta.Range("A10:B10").Copy()
$wsMeta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
